$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 13 ("Explicit Text:") so it becomes
# row 14, shifting every row below it down by one. The row that used to be
# row 12 (a blank spacer row) becomes the new row 13 spacer, and the newly
# inserted row 12 receives the new "Large Double Number:" data point.
$ws.Rows.Item(12).Insert()

$ws.Range("B12").Value = "Large Double Number:"
$ws.Range("C12").Value = 9.999 * [Math]::Pow(10, 307)
